# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, a new blank spacer column is inserted
# before the old "Late" column (N). This pushes the existing
# Late / heading / Outstanding columns one place to the right
# (N -> O, O -> P, P -> Q) and widens the sheet from A:P to A:Q.
# The sheet is also left as the active tab/sheet with K12 selected.

$wb = $excel.ActiveWorkbook

# "Repayment schedule" is the 3rd sheet (Input, Summary, Repayment schedule, Transactions).
$ws = $wb.Worksheets.Item(3)

# Insert a new blank column at N; everything from N onward shifts right.
$ws.Columns("N").Insert() | Out-Null

# Give the new spacer column the same width as column M (In Advance).
$ws.Columns("N").ColumnWidth = 9.83

# Leave the workbook with "Repayment schedule" as the active sheet/tab,
# with K12 as the selected cell.
$ws.Activate() | Out-Null
$ws.Range("K12").Select() | Out-Null
